$d = $word.ActiveDocument

# --- Change 1: paragraph "Classes are user-defined datatypes ..." ---
# Replace the trailing "." of "...class members." with an ellipsis "…"
# (Touching the Font property first keeps the existing run boundaries intact,
#  matching the target structure where "datatypes" stays its own run.)
$para1 = $d.Paragraphs(3)
$full1 = $para1.Range.Text
$len1 = $full1.Length
$periodStart = $para1.Range.Start + $len1 - 2
$periodEnd = $para1.Range.Start + $len1 - 1
$rngPeriod = $d.Range($periodStart, $periodEnd)
$rngPeriod.Font.Name = $rngPeriod.Font.Name
$rngPeriod.Text = "…"

# --- Change 2: paragraph "If one member function is called ... same class it is called nesting of a member function." ---
# Re-save the text spanning the "class" proofErr markers so the run/proofErr
# marks collapse into a single contiguous run (text itself is unchanged).
$d.Content.Find.Execute("same class it is called", $true, $false, $false, $false, $false, $true, 1, $false, "same class it is called", 2)

# --- Change 3: paragraph "The memory is only allocated ... So the objects don't have individual copies..." ---
# Re-save the text spanning the "So" proofErr markers so the run/proofErr
# marks collapse into a single contiguous run (text itself is unchanged).
$d.Content.Find.Execute("declared. So the objects", $true, $false, $false, $false, $false, $true, 1, $false, "declared. So the objects", 2)
